$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AGR")
$template = $wb.Worksheets.Item("UPD_biogas")

# Insert two new columns (Attrib_Cond, Val_Cond) before the region formula
# columns, then drop the old Cset_CN column (which shifts to the far right).
$ws.Range("F1:G1").EntireColumn.Insert()
$ws.Range("L1").EntireColumn.Delete()

# New header cells for the inserted columns.
$ws.Range("F3").Value = "Attrib_Cond"
$ws.Range("G3").Value = "Val_Cond"

# Rework the existing two rows: Attribute ACT_BND -> CAP_BND and
# Pset_PN ABIOFRSR* -> IMPBIOWPE_S2.
$ws.Range("D4").Value = "CAP_BND"
$ws.Range("D5").Value = "CAP_BND"
$ws.Range("K4").Value = "IMPBIOWPE_S2"
$ws.Range("K5").Value = "IMPBIOWPE_S2"

# Add three more pairs of rows disallowing imports above 2018 levels for
# the remaining biomass potentials.
$names = @("IMPBIOWPE_S3", "IMPBIOWCH_S2", "IMPBIOWCH_S3")
$row = 6
foreach ($name in $names) {
    $ws.Range("D$row").Value = "CAP_BND"
    $ws.Range("E$row").Value = 2018
    $ws.Range("H$row").Value = 0
    $ws.Range("I$row").Value = 0
    $ws.Range("K$row").Value = $name
    $row++

    $ws.Range("D$row").Value = "CAP_BND"
    $ws.Range("E$row").Value = 0
    $ws.Range("H$row").Value = 5
    $ws.Range("I$row").Value = 5
    $ws.Range("K$row").Value = $name
    $row++
}

# Add the i/e rule explanatory comment on I2, matching the format used on
# the UPD_biogas sheet (copy its formatting for I2:K2 onto this sheet).
$ws.Range("I2").AddComment("Define the qualifiers based upon technology set + topology + name + descriptions, according to both include and exclude specifications.")
$template.Range("I2:K2").Copy()
$ws.Range("I2").PasteSpecial(-4122)

# Widen column K to fit the longer Pset_PN qualifier names.
$ws.Columns.Item(11).AutoFit()

# Rename the sheet from AGR to IRE.
$ws.Name = "IRE"
